# gray-edges.pptx tweak
#
# Per the commit diff, slide 4 ("Gray Edges") has the entrance/fade
# animation on the *second* bullet of the "Content Placeholder 2" body
# text (paragraph range st=2/end=2, timing node id="15") slowed down
# from 500ms to 1000ms.
#
# (The diff also shows the VML "spid" attribute on several legacy
# Equation OLE objects -- p:oleObj spid="_x0000_s####" -- incrementing
# by a couple of numbers on slides 7, 11, 12, 13, 14, 15 and 22. That
# spid is an internal legacy-VML shape id minted by PowerPoint's own
# shape-id pool; it isn't part of the PowerPoint object model (no
# Shape/OLEFormat property reads or writes it, and it is not affected
# by moving/resizing/copying the shape), so it cannot be produced
# through COM automation and is intentionally left untouched here.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$seq = $s.TimeLine.MainSequence
# Effects on slide 4: 1) Group 4 wipe, 2) bullet-1 fade, 3) bullet-2 fade.
$effect = $seq.Item(3)
$effect.Timing.Duration = 1
